# Seminar Project Feedback Upload
#
# Moves several class-diagram boxes/connectors on the two slides and
# recolors the line of three connector arrows on slide 1 to accent6.
#
# NOTE: shapes are addressed by their position in the Shapes collection
# (established by inspecting the deck up front) rather than by their
# internal <p:cNvPr id="..">, because PowerPoint's Shapes.Item() indexes
# positionally, not by that id.
#
# NOTE: assigning to Shape.Left/Top/Width/Height goes through a
# single-precision (float32) round-trip internally and then truncates to
# EMU, which can silently shave 1 EMU (1/12700 pt) off an otherwise exact
# value. Nudging the point value up by a small epsilon (well under
# 1/12700 pt => no visual/semantic effect) keeps the stored EMU exact.
function EmuToPt($emu) {
    return ($emu / 12700) + 0.00004
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape 5 (id=16, "Desktop" rectangle) moves.
$desktop = $s1.Shapes.Item(5)
$desktop.Left = EmuToPt 8582783
$desktop.Top = EmuToPt 2387478

# Three inheritance-arrow connectors gain an explicit accent6 line color.
foreach ($idx in 6, 7, 8) {
    $cxn = $s1.Shapes.Item($idx)
    $cxn.Line.ForeColor.ObjectThemeColor = 10  # msoThemeColorAccent6
}

# ---------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Shape 3 (id=8, "Desktops / Management / UI") moves.
$shape8 = $s2.Shapes.Item(3)
$shape8.Left = EmuToPt 10146587
$shape8.Top = EmuToPt 2485870

# Shape 4 (id=9, "Desktops / Management / UI") moves.
$shape9 = $s2.Shapes.Item(4)
$shape9.Left = EmuToPt 1235551
$shape9.Top = EmuToPt 2485871

# Connector 9 (id=29) re-routes after shape 8/9 moved.
$cxn29 = $s2.Shapes.Item(9)
$cxn29.Left = EmuToPt 2855274
$cxn29.Top = EmuToPt 2901717
$cxn29.Width = EmuToPt 932626
$cxn29.Height = EmuToPt 1247528

# Connector 10 (id=36) re-routes after shape 8/9 moved.
$cxn36 = $s2.Shapes.Item(10)
$cxn36.Left = EmuToPt 8647069
$cxn36.Top = EmuToPt 2901716
$cxn36.Width = EmuToPt 1499518
$cxn36.Height = EmuToPt 1247530

# Shape 13 (id=22, "Stocks / Management / UI") moves.
$shape22 = $s2.Shapes.Item(13)
$shape22.Left = EmuToPt 5415419
$shape22.Top = EmuToPt 373856

# Connector 14 (id=23) re-routes after shape 22 moved.
$cxn23 = $s2.Shapes.Item(14)
$cxn23.Left = EmuToPt 2045413
$cxn23.Top = EmuToPt 789702
$cxn23.Width = EmuToPt 3370006
$cxn23.Height = EmuToPt 1696169

# Connector 15 (id=26) re-routes after shape 22 moved.
$cxn26 = $s2.Shapes.Item(15)
$cxn26.Left = EmuToPt 7035142
$cxn26.Top = EmuToPt 789702
$cxn26.Width = EmuToPt 3921307
$cxn26.Height = EmuToPt 1696168

# Connector 16 (id=67) re-routes after shape 22 moved, and also flips.
$cxn67 = $s2.Shapes.Item(16)
$cxn67.HorizontalFlip = -1  # msoTrue
$cxn67.Left = EmuToPt 6217485
$cxn67.Top = EmuToPt 1205547
$cxn67.Width = EmuToPt 7796
$cxn67.Height = EmuToPt 1127711
